$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2025-11-02 Sunday" "2025-11-03 Monday"

Replace-Text "814×6=4884" "530×8=4240"
Replace-Text "118×2=236" "249×4=996"
Replace-Text "843×2=1686" "477×6=2862"
Replace-Text "809×6=4854" "788×5=3940"
Replace-Text "819×4=3276" "988×5=4940"

Replace-Text "367×6=2202" "546×6=3276"
Replace-Text "358×6=2148" "120×3=360"
Replace-Text "448×2=896" "551×2=1102"
Replace-Text "600×4=2400" "634×5=3170"
Replace-Text "104×7=728" "996×5=4980"

Replace-Text "710×8=5680" "648×9=5832"
Replace-Text "950×2=1900" "631×4=2524"
Replace-Text "692×9=6228" "597×7=4179"
Replace-Text "719×8=5752" "974×7=6818"
Replace-Text "432×8=3456" "580×8=4640"

Replace-Text "738×8=5904" "175×2=350"
Replace-Text "390×7=2730" "250×6=1500"
Replace-Text "378×4=1512" "576×8=4608"
Replace-Text "896×9=8064" "141×3=423"
Replace-Text "559×3=1677" "390×5=1950"

Replace-Text "762×4=3048" "863×8=6904"
Replace-Text "669×7=4683" "672×2=1344"
Replace-Text "162×8=1296" "813×5=4065"
Replace-Text "548×3=1644" "415×8=3320"
Replace-Text "971×3=2913" "378×7=2646"
